# complex_validate_test.xlsx — "assign" now uses a 'calculation' column
# for the formula assigned to the name (the former 'default' column is
# removed from the survey/section1/section2 sheets). The now-unused
# "default" column header is dropped by deleting the whole column, which
# shifts every later column (choice_list_name/validation_tags/
# hideInContents/boolean flag columns, etc.) one slot to the left.

$wb = $excel.ActiveWorkbook

$survey = $wb.Worksheets.Item("survey")
$section1 = $wb.Worksheets.Item("section1")
$section2 = $wb.Worksheets.Item("section2")

# Remove the (empty, unused) "default" column -- column R on every sheet
# that has it -- shifting the trailing columns left by one.
$survey.Columns("R").Delete()
$section1.Columns("R").Delete()
$section2.Columns("R").Delete()

# Update the remembered selections on the edited sheets to match where the
# editor ended up after the column removal.
$survey.Range("R1:R1048576").Select()
$section1.Range("R1:R1048576").Select()
$section2.Range("P10").Select()

# "section2" (sheet index 3 / tab 2) ends up the active tab.
$section2.Activate()
